$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.638.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.170.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.08%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.609"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.41"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.82%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.580"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0909"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.75%  "
$ws.Range("E13").Value = "  -3.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.494.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.165.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.783"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.510.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.74%  "
$ws.Range("E20").Value = "  -1.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -12.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "227.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.50%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.94%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.88%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0772"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -10.42%  "
$ws.Range("E36").Value = "  -3.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.105"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0308"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.02"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -11.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "59.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.58%  "
$ws.Range("E45").Value = "  -5.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0964"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.60%  "
$ws.Range("E48").Value = "  -4.12%  "
$ws.Range("E49").Value = "  -5.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.39%  "
$ws.Range("E51").Value = "  -2.27%  "
